# Daily Report update: append the next day's (2026-02-13, serial 46066)
# depository rows to the Daily_Data sheet by repeating the most recent
# day's block of 22 rows (one Registered/Eligible pair per depository),
# advancing the date by one day.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Daily_Data")

# The most recently reported day occupies rows 618-639 (22 rows: 11
# depositories x Registered/Eligible). Duplicate that block into the new
# rows 640-661, incrementing the date by one day and copying every other
# column's value unchanged.
$srcStart = 618
$dstStart = 640
$rowCount = 22
$lastCol = 8

for ($i = 0; $i -lt $rowCount; $i++) {
    $srcRow = $srcStart + $i
    $dstRow = $dstStart + $i

    # Column A holds the report date (serial date number) - advance by 1 day.
    $srcDate = $ws.Cells.Item($srcRow, 1).Value2
    $ws.Cells.Item($dstRow, 1).Value = $srcDate + 1
    $ws.Cells.Item($dstRow, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"

    # Columns B (Region_Type label) through H (TOTAL_TODAY) are copied as-is.
    for ($col = 2; $col -le $lastCol; $col++) {
        $ws.Cells.Item($dstRow, $col).Value = $ws.Cells.Item($srcRow, $col).Value2
    }
}
